$d = $word.ActiveDocument

$d.Content.Find.Execute("724÷8=90, 4", $true, $false, $false, $false, $false, $true, 1, $false, "308÷8=38, 4", 2) | Out-Null
$d.Content.Find.Execute("356÷4=89, 0", $true, $false, $false, $false, $false, $true, 1, $false, "853÷6=142, 1", 2) | Out-Null
$d.Content.Find.Execute("442÷4=110, 2", $true, $false, $false, $false, $false, $true, 1, $false, "170÷6=28, 2", 2) | Out-Null
$d.Content.Find.Execute("898÷9=99, 7", $true, $false, $false, $false, $false, $true, 1, $false, "989÷4=247, 1", 2) | Out-Null
$d.Content.Find.Execute("596÷3=198, 2", $true, $false, $false, $false, $false, $true, 1, $false, "183÷9=20, 3", 2) | Out-Null
$d.Content.Find.Execute("568÷4=142, 0", $true, $false, $false, $false, $false, $true, 1, $false, "958÷3=319, 1", 2) | Out-Null
$d.Content.Find.Execute("238÷5=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "560÷6=93, 2", 2) | Out-Null
$d.Content.Find.Execute("273÷6=45, 3", $true, $false, $false, $false, $false, $true, 1, $false, "628÷2=314, 0", 2) | Out-Null
$d.Content.Find.Execute("497÷9=55, 2", $true, $false, $false, $false, $false, $true, 1, $false, "232÷3=77, 1", 2) | Out-Null
$d.Content.Find.Execute("959÷7=137, 0", $true, $false, $false, $false, $false, $true, 1, $false, "162÷4=40, 2", 2) | Out-Null
$d.Content.Find.Execute("387÷5=77, 2", $true, $false, $false, $false, $false, $true, 1, $false, "896÷3=298, 2", 2) | Out-Null
$d.Content.Find.Execute("300÷7=42, 6", $true, $false, $false, $false, $false, $true, 1, $false, "288÷5=57, 3", 2) | Out-Null
$d.Content.Find.Execute("506÷9=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "905÷7=129, 2", 2) | Out-Null
$d.Content.Find.Execute("610÷4=152, 2", $true, $false, $false, $false, $false, $true, 1, $false, "527÷9=58, 5", 2) | Out-Null
$d.Content.Find.Execute("585÷2=292, 1", $true, $false, $false, $false, $false, $true, 1, $false, "316÷7=45, 1", 2) | Out-Null
$d.Content.Find.Execute("946÷5=189, 1", $true, $false, $false, $false, $false, $true, 1, $false, "616÷4=154, 0", 2) | Out-Null
$d.Content.Find.Execute("746÷9=82, 8", $true, $false, $false, $false, $false, $true, 1, $false, "722÷6=120, 2", 2) | Out-Null
$d.Content.Find.Execute("682÷6=113, 4", $true, $false, $false, $false, $false, $true, 1, $false, "545÷4=136, 1", 2) | Out-Null
$d.Content.Find.Execute("359÷7=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "896÷7=128, 0", 2) | Out-Null
$d.Content.Find.Execute("768÷4=192, 0", $true, $false, $false, $false, $false, $true, 1, $false, "360÷7=51, 3", 2) | Out-Null
$d.Content.Find.Execute("938÷6=156, 2", $true, $false, $false, $false, $false, $true, 1, $false, "480÷8=60, 0", 2) | Out-Null
$d.Content.Find.Execute("141÷2=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "591÷5=118, 1", 2) | Out-Null
$d.Content.Find.Execute("774÷6=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "204÷4=51, 0", 2) | Out-Null
$d.Content.Find.Execute("407÷7=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "364÷3=121, 1", 2) | Out-Null
$d.Content.Find.Execute("835÷8=104, 3", $true, $false, $false, $false, $false, $true, 1, $false, "893÷8=111, 5", 2) | Out-Null
